$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 3910.501516717734
$ws.Range("D2").Value = 570.9193558457891
